$wb = $excel.ActiveWorkbook

# --- Sheet "Inputs": add a new "Active" column (D) ---
$ws = $wb.Worksheets.Item("Inputs")

$ws.Range("D1").Value = "Active"
$ws.Range("D2").Value = $true

# --- Sheet "Notes": tweak the generated-fixture note text ---
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A1").Value = "Generated from fixture workbook"
